$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content ---

# Row 10: "Objetivos:" row - B/C text changes to the docente name (string reused elsewhere)
$ws.Range("B10").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C10").Value = "6495737 - Durval Rodrigues Junior"

# Row 13: becomes "Programa resumido:" in A (new), B/C keep being filled (previously held first docente name)
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "6495737 - Durval Rodrigues Junior"
$ws.Range("C13").Value = "6495737 - Durval Rodrigues Junior"

# Row 14: becomes "Short syllabus:" (A only now; drop the previous B/C docente name entirely)
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# Row 15: becomes "Programa:" in A (new), B/C become the second docente name
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "5983729 - Fernando Vernilli Junior"
$ws.Range("C15").Value = "5983729 - Fernando Vernilli Junior"

# Row 16: becomes "Syllabus:" (A only now; drop the previous B/C content)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# Row 17: becomes "Avaliação:" (A only, already A-only before)
$ws.Range("A17").Value = "Avaliação:"

# Row 18: becomes "Método:" / third docente name (A/B/C all pre-existing)
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C18").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# Row 19: becomes "Critério:" (A pre-existing) / exam statement (B/C are NEW cells)
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."
$ws.Range("C19").Value = "Serão aplicadas duas provas escritas com notas P1 e P2."

# Row 20: becomes "Norma de recuperação:" (A pre-existing) / grade formula (B/C are NEW cells)
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."
$ws.Range("C20").Value = "A nota final NF será calculada pela fórmula: NF=(P1 + P2)/2."

# Row 21: becomes "Bibliografia:" / recovery rule text (A/B/C all pre-existing)
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."
$ws.Range("C21").Value = "Será aplicada uma prova escrita NR que comporá com a nota final NF a média final após recuperação MF=(NF+NF)/2."

# --- Fix formatting on brand-new cells: newly-created cells inherit the row's
# --- dominant style rather than their own column style, so re-apply the
# --- correct column format (copied from an existing, correctly-styled cell
# --- in the same column) without touching the values already written above.
$ws.Range("B10").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the trailing rows that no longer exist (old rows 22-24) ---
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

# --- Row heights (ht="..." customHeight="1"); row 17 reverts to default (no custom height) ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
